$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$ws.Activate()

# Daily sprint meeting for day 10 (column P):
# Task on row 12 ("Review code", To-do item) moved from "In progress" to "Done"
# and 1 hour of effort was logged against Day 10 for it.
$ws.Range("F12").Value = "Done"
$ws.Range("P12").Value = 1

# Reflect the new active cell/selection left by the edit (day-10 column, row 12)
$ws.Range("Q12").Select()

$excel.Calculate()
$excel.CalculateFull()

# Best-effort: nudge the burndown chart on the "Chart" sheet to refresh its
# cached series data now that Sprint!P12 / the Day-10 totals changed.
$cws = $wb.Worksheets.Item("Chart")
if ($cws.ChartObjects().Count() -gt 0) {
    $co = $cws.ChartObjects().Item(1)
    $co.Chart.Refresh()
}
